# Migrate the testing-framework label on the architecture diagram from
# "Karma" to "Jest" (switch from Karma/Jasmine to Jest for JS tests).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the textbox that currently reads "Karma" and update its text.
$found = $false
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        if ($shape.TextFrame.HasText) {
            if ($shape.TextFrame.TextRange.Text -eq "Karma") {
                $shape.TextFrame.TextRange.Text = "Jest"
                $found = $true
            }
        }
    }
}

if (-not $found) {
    Write-Host "Warning: 'Karma' textbox not found"
}
